$d = $word.ActiveDocument

# Create the three new "join" paragraph styles, each based on the
# existing MSCJoin style. MSC_Join_B additionally carries the CJK font
# run-properties (mirroring MSC_Paragraph_B).
$joinBase = $d.Styles("MSCJoin")

$sA = $d.Styles.Add("MSC_Join_A", 1)
$sA.BaseStyle = $joinBase

$sB = $d.Styles.Add("MSC_Join_B", 1)
$sB.BaseStyle = $joinBase
$sB.Font.Name = "Noto Sans CJK SC"
$sB.Font.NameFarEast = "Noto Sans CJK SC"
$sB.Font.NameBi = "Noto Sans CJK SC"

$sC = $d.Styles.Add("MSC_Join_C", 1)
$sC.BaseStyle = $joinBase

# Walk every paragraph in document order. Each "join" paragraph
# (style MSCJoin) sits inside the table cell that most recently started
# with an MSC_Paragraph_A / _B / _C paragraph - retarget it to the
# matching MSC_Join_A / _B / _C style.
$current = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "MSC_Paragraph_A") {
        $current = "A"
    } elseif ($styleName -eq "MSC_Paragraph_B") {
        $current = "B"
    } elseif ($styleName -eq "MSC_Paragraph_C") {
        $current = "C"
    } elseif ($styleName -eq "MSC_Join") {
        if ($current -eq "A") {
            $p.Style = $sA
        } elseif ($current -eq "B") {
            $p.Style = $sB
        } elseif ($current -eq "C") {
            $p.Style = $sC
        }
    }
}
